$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.5216437782290682
$ws.Range("D2").Value = 0.6071307011250271

# Row 3
$ws.Range("C3").Value = 0.003599592078948115
$ws.Range("D3").Value = 0.9971603886778233

# Row 4
$ws.Range("C4").Value = 0.8504632829603186
$ws.Range("D4").Value = 0.4042289819424913

# Row 5
$ws.Range("C5").Value = 4.230116814967821
$ws.Range("D5").Value = 0.0003439395751421159

# Row 6
$ws.Range("C6").Value = -0.2881821710824145
$ws.Range("D6").Value = 0.775904859031268

# Row 7
$ws.Range("C7").Value = 0.5057823489978439
$ws.Range("D7").Value = 0.6180398756601588

# Row 8
$ws.Range("C8").Value = 3.411135144805184
$ws.Range("D8").Value = 0.002503861860086998

# Row 9
$ws.Range("C9").Value = 0.7653533978657637
$ws.Range("D9").Value = 0.4521942631633329

# Row 10
$ws.Range("C10").Value = 1.759987207046143
$ws.Range("D10").Value = 0.09230901304066408
$ws.Range("G10").Value = "No"

# Row 11
$ws.Range("C11").Value = 2.083164072096575
$ws.Range("D11").Value = 0.04906760228396001
